# Remove anchor feature prioritization in feature selection
# Update Cox results table values for rows 11-16 (Proposed / Advanced sections)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "0.599 (nan, nan)"
$ws.Range("C11").Value = "0.567 (nan, nan)"
$ws.Range("D11").Value = "0.593 (nan, nan)"

$ws.Range("B12").Value = "0.589 (nan, nan)"
$ws.Range("C12").Value = "0.600 (nan, nan)"
$ws.Range("D12").Value = "0.533 (nan, nan)"

$ws.Range("B13").Value = "0.552 (nan, nan)"
$ws.Range("D13").Value = "0.563 (nan, nan)"

$ws.Range("B14").Value = "0.599 (nan, nan)"
$ws.Range("D14").Value = "0.607 (nan, nan)"

$ws.Range("B15").Value = "0.532 (nan, nan)"
$ws.Range("C15").Value = "0.667 (nan, nan)"

$ws.Range("B16").Value = "0.612 (nan, nan)"
$ws.Range("D16").Value = "0.533 (nan, nan)"
